$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update so_id (G), host_organization (H), issn_l (I) for rows 10-12
foreach ($row in 10..12) {
    $ws.Cells.Item($row, 7).Value = "https://openalex.org/S4386621754"
    $ws.Cells.Item($row, 8).Value = "Elsevier BV"
    $ws.Cells.Item($row, 9).Value = "2666-6367"
}
